# C5-PowerPoint.pptx — change the table style applied to the "SOURCES OF
# FINANCE" table (slide 6) from the custom default table style
# {024F96CC-0DD3-4055-884C-AD49224BEE36} to the built-in gallery style
# {7216D6F3-18C1-47EE-B77A-9B4097869697}.
#
# We search every slide/shape for the table that currently carries the old
# style id rather than hard-coding slide/shape indices, so the script is
# resilient to any reordering.

$p = $ppt.ActivePresentation

$oldStyleId = "{024F96CC-0DD3-4055-884C-AD49224BEE36}"
$newStyleId = "{7216D6F3-18C1-47EE-B77A-9B4097869697}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.HasTable) {
            $tbl = $shape.Table

            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
